# Weekly update: add a new price-report week (2023-07-25, Excel serial 45132)
# for "Fruta, Terminal La Palmera de La Serena - Plátano".
#
# The source data is laid out as repeating 3-row blocks (Pintón / Primera
# Maduro / Primera Pintón) ordered with the most recent week first. Adding
# the new week means inserting a fresh 3-row block right after the header,
# above the block that is currently first (row 1144), which pushes every
# existing block down by 3 rows (dimension grows from T1254 to T1257).
#
# The new block carries forward the same Volumen/Precio values that were
# in last week's (now second) block, only the Fecha changes to the new
# report date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 3 blank rows before the current first data block (row 1144).
#    Excel shifts every row at/after 1144 down by 3, so the old 1144:1146
#    block (and everything after it) now lives at 1147:1149, etc.
$ws.Rows("1144:1146").Insert()

# 2) Seed the new block by copying the (shifted) old first block's values
#    and formatting down into the freshly inserted rows.
$ws.Range("A1147:T1149").Copy($ws.Range("A1144:T1146"))

# 3) Stamp the new block with this week's report date.
$ws.Range("D1144:D1146").Value = 45132
